# "Generate Report for Handoff"
#
# Updates the localization-status report after a fresh handoff report run:
#   - Status moves from "Handed back: in sync with en-US" to "In Translation"
#     (both on the Overview rollup and on each language detail sheet).
#   - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" stamps
#     are refreshed.
#   - Each language sheet gets an Error Detail note about the handback file
#     being stale.
#   - A couple of columns that used to be sized for long timestamps are
#     narrowed back down (others are widened for the new Error Detail text).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$errorDetail = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0729dc005dfb2c635e2cf1a74b23e5cacd7ace06/e2e/ae7f0526-159b-4eaf-aafd-6e77a2be2935.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7299dbfcb75c8adee08f3f105ef0eab6c5e1c712/e2e/ae7f0526-159b-4eaf-aafd-6e77a2be2935.md.'

# ---------------------------------------------------------------------
# Overview sheet: per-language status + the shared "last generated" date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2017-02-09 14:17:41"

$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2017-02-09 14:17:22"
$wsZhCn.Range("R2").Value = $errorDetail

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsZhCn.Columns.Item(18).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2017-02-09 14:17:41"
$wsDeDe.Range("R2").Value = $errorDetail

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(18).ColumnWidth = 39.16666666666667
